$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B -> C, old C -> D),
# shifting the "Jun_13"/"Jun_10" rating columns to the right to make
# room for a new "Jun_27" rating column.
$ws.Columns("B").Insert()

# New column header + "UN" fill for the new Jun_27 column (rows 2-27
# mirror the existing UN placeholder values used throughout the sheet).
$ws.Range("B1").Value = "Jun_27"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# Resize columns: A stays the same, B/C (the two "UN" rating columns)
# get narrower, D (the widest, most-recent comments column) gets wider.
$ws.Columns("B").ColumnWidth = 26
$ws.Columns("C").ColumnWidth = 26
$ws.Columns("D").ColumnWidth = 51.1

# Add two new rows for the newly tracked firms.
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"

# Move the active selection as recorded by the author.
$ws.Range("B5").Select()
